$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Default the UserStatus column (N) to an empty string for all data rows (2-100)
$ws.Range("N2:N100").Value = ""

# Update the sheet's active selection to N1
$ws.Range("N1").Select()
